$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '68.316.26'
Set-TextValue $ws.Range('E2') '  -1.00%  '
Set-TextValue $ws.Range('D3') '3.900.22'
Set-TextValue $ws.Range('E3') '  -0.66%  '
Set-TextValue $ws.Range('E4') '  -0.09%  '
Set-TextValue $ws.Range('D5') '486.53'
Set-TextValue $ws.Range('E5') '  +0.59%  '
Set-TextValue $ws.Range('D6') '145.75'
Set-TextValue $ws.Range('E6') '  -0.61%  '
Set-TextValue $ws.Range('E7') '  +0.10%  '
Set-TextValue $ws.Range('D8') '0.998'
Set-TextValue $ws.Range('E8') '  +0.03%  '
Set-TextValue $ws.Range('D9') '0.742'
Set-TextValue $ws.Range('E9') '  +2.70%  '
Set-TextValue $ws.Range('D10') '0.180'
Set-TextValue $ws.Range('E10') '  +5.56%  '
Set-TextValue $ws.Range('D11') '0.0000355'
Set-TextValue $ws.Range('E11') '  -1.40%  '
Set-TextValue $ws.Range('D12') '42.94'
Set-TextValue $ws.Range('E12') '  +0.79%  '
Set-TextValue $ws.Range('D13') '10.50'
Set-TextValue $ws.Range('E13') '  +0.56%  '
Set-TextValue $ws.Range('D14') '4.519.93'
Set-TextValue $ws.Range('E14') '  -0.66%  '
Set-TextValue $ws.Range('D15') '3.903.71'
Set-TextValue $ws.Range('E15') '  -0.68%  '
Set-TextValue $ws.Range('D16') '14.23'
Set-TextValue $ws.Range('E16') '  -2.31%  '
Set-TextValue $ws.Range('E17') '  -0.59%  '
Set-TextValue $ws.Range('D18') '20.01'
Set-TextValue $ws.Range('E18') '  +1.71%  '
Set-TextValue $ws.Range('E19') '  +0.74%  '
Set-TextValue $ws.Range('D20') '68.327.36'
Set-TextValue $ws.Range('E20') '  -1.08%  '
Set-TextValue $ws.Range('D21') '431.58'
Set-TextValue $ws.Range('E21') '  -0.62%  '
Set-TextValue $ws.Range('E22') '  +6.71%  '
Set-TextValue $ws.Range('D23') '14.75'
Set-TextValue $ws.Range('E23') '  +1.04%  '
Set-TextValue $ws.Range('D24') '12.39'
Set-TextValue $ws.Range('E24') '  +19.72%  '
Set-TextValue $ws.Range('D25') '89.05'
Set-TextValue $ws.Range('E25') '  +1.23%  '
Set-TextValue $ws.Range('D26') '3.72'
Set-TextValue $ws.Range('E26') '  +4.40%  '
Set-TextValue $ws.Range('D27') '11.00'
Set-TextValue $ws.Range('E27') '  -5.57%  '
Set-TextValue $ws.Range('D28') '37.31'
Set-TextValue $ws.Range('E28') '  -2.66%  '
Set-TextValue $ws.Range('D29') '5.68'
Set-TextValue $ws.Range('E29') '  -3.96%  '
Set-TextValue $ws.Range('D30') '720.40'
Set-TextValue $ws.Range('E30') '  +1.60%  '
Set-TextValue $ws.Range('E31') '  +1.69%  '
Set-TextValue $ws.Range('E32') '  +0.73%  '
Set-TextValue $ws.Range('D33') '2.93'
Set-TextValue $ws.Range('E33') '  +2.83%  '
Set-TextValue $ws.Range('D34') '61.69'
Set-TextValue $ws.Range('E34') '  +5.08%  '
Set-TextValue $ws.Range('D35') '0.0₃0879'
Set-TextValue $ws.Range('E35') '  -6.65%  '
Set-TextValue $ws.Range('D36') '6.06'
Set-TextValue $ws.Range('E36') '  +8.71%  '
Set-TextValue $ws.Range('E37') '  -0.91%  '
Set-TextValue $ws.Range('D38') '0.400'
Set-TextValue $ws.Range('E38') '  +17.11%  '
Set-TextValue $ws.Range('E39') '  -3.64%  '
Set-TextValue $ws.Range('E40') '  +0.02%  '
Set-TextValue $ws.Range('D41') '0.0497'
Set-TextValue $ws.Range('E41') '  +5.47%  '
Set-TextValue $ws.Range('E42') '  +7.38%  '
Set-TextValue $ws.Range('D43') '3.08'
Set-TextValue $ws.Range('E43') '  +2.82%  '
Set-TextValue $ws.Range('D44') '3.00'
Set-TextValue $ws.Range('E44') '  -0.33%  '
Set-TextValue $ws.Range('B45') 'Stellar'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D45') '0.142'
Set-TextValue $ws.Range('E45') '  +1.10%  '
Set-TextValue $ws.Range('B46') 'BabyDogeCoin'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D46') '0.0₆0369'
Set-TextValue $ws.Range('E46') '  +33.00%  '
Set-TextValue $ws.Range('E47') '  -0.13%  '
Set-TextValue $ws.Range('D48') '3.33'
Set-TextValue $ws.Range('E48') '  +6.72%  '
Set-TextValue $ws.Range('E49') '  -1.12%  '
Set-TextValue $ws.Range('D50') '2.09'
Set-TextValue $ws.Range('E50') '  -2.54%  '
Set-TextValue $ws.Range('D51') '144.28'
Set-TextValue $ws.Range('E51') '  -2.54%  '
